# Regenerate the "within100" arithmetic-drill table: every equation cell's
# text is replaced with a newly generated one. Each Find/Replace targets the
# exact "<old>=" string of a single cell (all 100 old equation strings are
# unique in the document, so MatchWholeWord=$false Find is safe) and runs
# against the whole document content each time.
#
# NOTE: two calls below are intentionally reordered relative to the on-page
# top-to-bottom cell order so that no find-string is ever a substring of a
# not-yet-processed find-string (e.g. "6+61=" is contained in "16+61=") -
# that ensures Find can't accidentally match inside a cell meant for a later
# step.
$d = $word.ActiveDocument
$d.Content.Find.Execute("6+40=", $true, $false, $false, $false, $false, $true, 1, $false, "53+37=", 2) | Out-Null
$d.Content.Find.Execute("33+55=", $true, $false, $false, $false, $false, $true, 1, $false, "58+16=", 2) | Out-Null
$d.Content.Find.Execute("22+18=", $true, $false, $false, $false, $false, $true, 1, $false, "69-50=", 2) | Out-Null
$d.Content.Find.Execute("60-14=", $true, $false, $false, $false, $false, $true, 1, $false, "62-34=", 2) | Out-Null
$d.Content.Find.Execute("16+37=", $true, $false, $false, $false, $false, $true, 1, $false, "89-74=", 2) | Out-Null
$d.Content.Find.Execute("34-31=", $true, $false, $false, $false, $false, $true, 1, $false, "29+24=", 2) | Out-Null
$d.Content.Find.Execute("80-62=", $true, $false, $false, $false, $false, $true, 1, $false, "57+2=", 2) | Out-Null
$d.Content.Find.Execute("50+18=", $true, $false, $false, $false, $false, $true, 1, $false, "9+12=", 2) | Out-Null
$d.Content.Find.Execute("50+8=", $true, $false, $false, $false, $false, $true, 1, $false, "36+13=", 2) | Out-Null
$d.Content.Find.Execute("31-15=", $true, $false, $false, $false, $false, $true, 1, $false, "54+36=", 2) | Out-Null
$d.Content.Find.Execute("87+5=", $true, $false, $false, $false, $false, $true, 1, $false, "1+56=", 2) | Out-Null
$d.Content.Find.Execute("35+51=", $true, $false, $false, $false, $false, $true, 1, $false, "66-45=", 2) | Out-Null
$d.Content.Find.Execute("44-18=", $true, $false, $false, $false, $false, $true, 1, $false, "71-16=", 2) | Out-Null
$d.Content.Find.Execute("14+12=", $true, $false, $false, $false, $false, $true, 1, $false, "48-24=", 2) | Out-Null
$d.Content.Find.Execute("4+93=", $true, $false, $false, $false, $false, $true, 1, $false, "45+26=", 2) | Out-Null
$d.Content.Find.Execute("92-22=", $true, $false, $false, $false, $false, $true, 1, $false, "52-31=", 2) | Out-Null
$d.Content.Find.Execute("36-10=", $true, $false, $false, $false, $false, $true, 1, $false, "6+3=", 2) | Out-Null
$d.Content.Find.Execute("93-67=", $true, $false, $false, $false, $false, $true, 1, $false, "66-29=", 2) | Out-Null
$d.Content.Find.Execute("31+65=", $true, $false, $false, $false, $false, $true, 1, $false, "70-10=", 2) | Out-Null
$d.Content.Find.Execute("66-6=", $true, $false, $false, $false, $false, $true, 1, $false, "33+24=", 2) | Out-Null
$d.Content.Find.Execute("47-33=", $true, $false, $false, $false, $false, $true, 1, $false, "37+46=", 2) | Out-Null
$d.Content.Find.Execute("80-14=", $true, $false, $false, $false, $false, $true, 1, $false, "72-41=", 2) | Out-Null
$d.Content.Find.Execute("81-35=", $true, $false, $false, $false, $false, $true, 1, $false, "41+49=", 2) | Out-Null
$d.Content.Find.Execute("44-42=", $true, $false, $false, $false, $false, $true, 1, $false, "31-9=", 2) | Out-Null
$d.Content.Find.Execute("59-30=", $true, $false, $false, $false, $false, $true, 1, $false, "71+24=", 2) | Out-Null
$d.Content.Find.Execute("88-5=", $true, $false, $false, $false, $false, $true, 1, $false, "29+63=", 2) | Out-Null
$d.Content.Find.Execute("37+19=", $true, $false, $false, $false, $false, $true, 1, $false, "96-78=", 2) | Out-Null
$d.Content.Find.Execute("82-71=", $true, $false, $false, $false, $false, $true, 1, $false, "10+3=", 2) | Out-Null
$d.Content.Find.Execute("19+28=", $true, $false, $false, $false, $false, $true, 1, $false, "59-5=", 2) | Out-Null
$d.Content.Find.Execute("68-30=", $true, $false, $false, $false, $false, $true, 1, $false, "97-15=", 2) | Out-Null
$d.Content.Find.Execute("27-5=", $true, $false, $false, $false, $false, $true, 1, $false, "67+13=", 2) | Out-Null
$d.Content.Find.Execute("49-42=", $true, $false, $false, $false, $false, $true, 1, $false, "95-83=", 2) | Out-Null
$d.Content.Find.Execute("58-56=", $true, $false, $false, $false, $false, $true, 1, $false, "68-50=", 2) | Out-Null
$d.Content.Find.Execute("47+47=", $true, $false, $false, $false, $false, $true, 1, $false, "10+57=", 2) | Out-Null
$d.Content.Find.Execute("47-40=", $true, $false, $false, $false, $false, $true, 1, $false, "4+13=", 2) | Out-Null
$d.Content.Find.Execute("40+21=", $true, $false, $false, $false, $false, $true, 1, $false, "63-53=", 2) | Out-Null
$d.Content.Find.Execute("7-1=", $true, $false, $false, $false, $false, $true, 1, $false, "62-53=", 2) | Out-Null
$d.Content.Find.Execute("10+82=", $true, $false, $false, $false, $false, $true, 1, $false, "38+21=", 2) | Out-Null
$d.Content.Find.Execute("48+50=", $true, $false, $false, $false, $false, $true, 1, $false, "72-47=", 2) | Out-Null
$d.Content.Find.Execute("68-24=", $true, $false, $false, $false, $false, $true, 1, $false, "4+43=", 2) | Out-Null
$d.Content.Find.Execute("81-5=", $true, $false, $false, $false, $false, $true, 1, $false, "73-12=", 2) | Out-Null
$d.Content.Find.Execute("56+17=", $true, $false, $false, $false, $false, $true, 1, $false, "56+31=", 2) | Out-Null
$d.Content.Find.Execute("45-40=", $true, $false, $false, $false, $false, $true, 1, $false, "49-0=", 2) | Out-Null
$d.Content.Find.Execute("16+61=", $true, $false, $false, $false, $false, $true, 1, $false, "44-41=", 2) | Out-Null
$d.Content.Find.Execute("6+61=", $true, $false, $false, $false, $false, $true, 1, $false, "14+69=", 2) | Out-Null
$d.Content.Find.Execute("85+4=", $true, $false, $false, $false, $false, $true, 1, $false, "73+18=", 2) | Out-Null
$d.Content.Find.Execute("71-69=", $true, $false, $false, $false, $false, $true, 1, $false, "40+51=", 2) | Out-Null
$d.Content.Find.Execute("37+32=", $true, $false, $false, $false, $false, $true, 1, $false, "89-53=", 2) | Out-Null
$d.Content.Find.Execute("55-38=", $true, $false, $false, $false, $false, $true, 1, $false, "34+35=", 2) | Out-Null
$d.Content.Find.Execute("89-61=", $true, $false, $false, $false, $false, $true, 1, $false, "20-19=", 2) | Out-Null
$d.Content.Find.Execute("49-4=", $true, $false, $false, $false, $false, $true, 1, $false, "40+47=", 2) | Out-Null
$d.Content.Find.Execute("75-8=", $true, $false, $false, $false, $false, $true, 1, $false, "71-5=", 2) | Out-Null
$d.Content.Find.Execute("76-48=", $true, $false, $false, $false, $false, $true, 1, $false, "14+49=", 2) | Out-Null
$d.Content.Find.Execute("45-24=", $true, $false, $false, $false, $false, $true, 1, $false, "39+9=", 2) | Out-Null
$d.Content.Find.Execute("0+24=", $true, $false, $false, $false, $false, $true, 1, $false, "64-29=", 2) | Out-Null
$d.Content.Find.Execute("9+20=", $true, $false, $false, $false, $false, $true, 1, $false, "61-17=", 2) | Out-Null
$d.Content.Find.Execute("68-9=", $true, $false, $false, $false, $false, $true, 1, $false, "9+68=", 2) | Out-Null
$d.Content.Find.Execute("26+20=", $true, $false, $false, $false, $false, $true, 1, $false, "16+12=", 2) | Out-Null
$d.Content.Find.Execute("86-17=", $true, $false, $false, $false, $false, $true, 1, $false, "79+14=", 2) | Out-Null
$d.Content.Find.Execute("40+37=", $true, $false, $false, $false, $false, $true, 1, $false, "49-33=", 2) | Out-Null
$d.Content.Find.Execute("88-24=", $true, $false, $false, $false, $false, $true, 1, $false, "41+55=", 2) | Out-Null
$d.Content.Find.Execute("72+2=", $true, $false, $false, $false, $false, $true, 1, $false, "66+24=", 2) | Out-Null
$d.Content.Find.Execute("67-24=", $true, $false, $false, $false, $false, $true, 1, $false, "81-67=", 2) | Out-Null
$d.Content.Find.Execute("76-13=", $true, $false, $false, $false, $false, $true, 1, $false, "42-33=", 2) | Out-Null
$d.Content.Find.Execute("26-23=", $true, $false, $false, $false, $false, $true, 1, $false, "63-48=", 2) | Out-Null
$d.Content.Find.Execute("5-3=", $true, $false, $false, $false, $false, $true, 1, $false, "39+20=", 2) | Out-Null
$d.Content.Find.Execute("30+36=", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=", 2) | Out-Null
$d.Content.Find.Execute("32-29=", $true, $false, $false, $false, $false, $true, 1, $false, "31+60=", 2) | Out-Null
$d.Content.Find.Execute("94-47=", $true, $false, $false, $false, $false, $true, 1, $false, "85-22=", 2) | Out-Null
$d.Content.Find.Execute("85-71=", $true, $false, $false, $false, $false, $true, 1, $false, "45+29=", 2) | Out-Null
$d.Content.Find.Execute("9+39=", $true, $false, $false, $false, $false, $true, 1, $false, "73+2=", 2) | Out-Null
$d.Content.Find.Execute("84-36=", $true, $false, $false, $false, $false, $true, 1, $false, "40-18=", 2) | Out-Null
$d.Content.Find.Execute("28+56=", $true, $false, $false, $false, $false, $true, 1, $false, "54+36=", 2) | Out-Null
$d.Content.Find.Execute("63-15=", $true, $false, $false, $false, $false, $true, 1, $false, "37-31=", 2) | Out-Null
$d.Content.Find.Execute("14+3=", $true, $false, $false, $false, $false, $true, 1, $false, "74-26=", 2) | Out-Null
$d.Content.Find.Execute("2+97=", $true, $false, $false, $false, $false, $true, 1, $false, "43+47=", 2) | Out-Null
$d.Content.Find.Execute("40-19=", $true, $false, $false, $false, $false, $true, 1, $false, "64+33=", 2) | Out-Null
$d.Content.Find.Execute("65+27=", $true, $false, $false, $false, $false, $true, 1, $false, "94-7=", 2) | Out-Null
$d.Content.Find.Execute("31+44=", $true, $false, $false, $false, $false, $true, 1, $false, "1+76=", 2) | Out-Null
$d.Content.Find.Execute("9+23=", $true, $false, $false, $false, $false, $true, 1, $false, "42-14=", 2) | Out-Null
$d.Content.Find.Execute("52+19=", $true, $false, $false, $false, $false, $true, 1, $false, "0+78=", 2) | Out-Null
$d.Content.Find.Execute("3+13=", $true, $false, $false, $false, $false, $true, 1, $false, "3+92=", 2) | Out-Null
$d.Content.Find.Execute("98-74=", $true, $false, $false, $false, $false, $true, 1, $false, "53-39=", 2) | Out-Null
$d.Content.Find.Execute("17+51=", $true, $false, $false, $false, $false, $true, 1, $false, "2+28=", 2) | Out-Null
$d.Content.Find.Execute("82-11=", $true, $false, $false, $false, $false, $true, 1, $false, "50-38=", 2) | Out-Null
$d.Content.Find.Execute("78-3=", $true, $false, $false, $false, $false, $true, 1, $false, "53-7=", 2) | Out-Null
$d.Content.Find.Execute("28+14=", $true, $false, $false, $false, $false, $true, 1, $false, "35+59=", 2) | Out-Null
$d.Content.Find.Execute("87-30=", $true, $false, $false, $false, $false, $true, 1, $false, "90-35=", 2) | Out-Null
$d.Content.Find.Execute("38+6=", $true, $false, $false, $false, $false, $true, 1, $false, "8+41=", 2) | Out-Null
$d.Content.Find.Execute("46-36=", $true, $false, $false, $false, $false, $true, 1, $false, "54+20=", 2) | Out-Null
$d.Content.Find.Execute("21+46=", $true, $false, $false, $false, $false, $true, 1, $false, "70-26=", 2) | Out-Null
$d.Content.Find.Execute("0+87=", $true, $false, $false, $false, $false, $true, 1, $false, "96-35=", 2) | Out-Null
$d.Content.Find.Execute("30+28=", $true, $false, $false, $false, $false, $true, 1, $false, "36-33=", 2) | Out-Null
$d.Content.Find.Execute("97-46=", $true, $false, $false, $false, $false, $true, 1, $false, "83-76=", 2) | Out-Null
$d.Content.Find.Execute("84-43=", $true, $false, $false, $false, $false, $true, 1, $false, "33+59=", 2) | Out-Null
$d.Content.Find.Execute("12-7=", $true, $false, $false, $false, $false, $true, 1, $false, "65-63=", 2) | Out-Null
$d.Content.Find.Execute("21+72=", $true, $false, $false, $false, $false, $true, 1, $false, "37-11=", 2) | Out-Null
$d.Content.Find.Execute("76-29=", $true, $false, $false, $false, $false, $true, 1, $false, "16+70=", 2) | Out-Null
$d.Content.Find.Execute("86-35=", $true, $false, $false, $false, $false, $true, 1, $false, "44+36=", 2) | Out-Null
$d.Content.Find.Execute("73-18=", $true, $false, $false, $false, $false, $true, 1, $false, "48-28=", 2) | Out-Null
